$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 corresponds to Year 2021 / Country Spain.
# Update Excess and its credible interval, and Percent_excess and its credible interval
# with the new results. The "Excess" and "Percent_excess" values look like plain numbers
# / percentages to Excel, so force them to be stored as text (matching the original
# shared-string table, which holds these as literal text, e.g. "10,836") before writing
# the new values - otherwise Excel would silently convert them to numeric cells.
$ws.Range("C15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"

$ws.Range("C15").Value = "10,867"
$ws.Range("D15").Value = "(-3,423; 24,715)"
$ws.Range("G15").Value = "5.0%"
$ws.Range("H15").Value = "(-1.4%; 11.8%)"
